$d = $word.ActiveDocument

$d.Content.Find.Execute("911×2=1822", $true, $false, $false, $false, $false, $true, 1, $false, "726×7=5082", 2) | Out-Null
$d.Content.Find.Execute("559×8=4472", $true, $false, $false, $false, $false, $true, 1, $false, "651×6=3906", 2) | Out-Null
$d.Content.Find.Execute("861×8=6888", $true, $false, $false, $false, $false, $true, 1, $false, "929×4=3716", 2) | Out-Null
$d.Content.Find.Execute("812×9=7308", $true, $false, $false, $false, $false, $true, 1, $false, "295×4=1180", 2) | Out-Null
$d.Content.Find.Execute("622×7=4354", $true, $false, $false, $false, $false, $true, 1, $false, "645×7=4515", 2) | Out-Null
$d.Content.Find.Execute("942×5=4710", $true, $false, $false, $false, $false, $true, 1, $false, "785×8=6280", 2) | Out-Null
$d.Content.Find.Execute("808×2=1616", $true, $false, $false, $false, $false, $true, 1, $false, "980×9=8820", 2) | Out-Null
$d.Content.Find.Execute("868×8=6944", $true, $false, $false, $false, $false, $true, 1, $false, "764×7=5348", 2) | Out-Null
$d.Content.Find.Execute("925×4=3700", $true, $false, $false, $false, $false, $true, 1, $false, "733×2=1466", 2) | Out-Null
$d.Content.Find.Execute("500×2=1000", $true, $false, $false, $false, $false, $true, 1, $false, "112×4=448", 2) | Out-Null
$d.Content.Find.Execute("211×5=1055", $true, $false, $false, $false, $false, $true, 1, $false, "948×3=2844", 2) | Out-Null
$d.Content.Find.Execute("828×8=6624", $true, $false, $false, $false, $false, $true, 1, $false, "736×9=6624", 2) | Out-Null
$d.Content.Find.Execute("216×4=864", $true, $false, $false, $false, $false, $true, 1, $false, "627×8=5016", 2) | Out-Null
$d.Content.Find.Execute("736×2=1472", $true, $false, $false, $false, $false, $true, 1, $false, "645×8=5160", 2) | Out-Null
$d.Content.Find.Execute("949×3=2847", $true, $false, $false, $false, $false, $true, 1, $false, "431×6=2586", 2) | Out-Null
$d.Content.Find.Execute("743×8=5944", $true, $false, $false, $false, $false, $true, 1, $false, "948×8=7584", 2) | Out-Null
$d.Content.Find.Execute("880×8=7040", $true, $false, $false, $false, $false, $true, 1, $false, "390×4=1560", 2) | Out-Null
$d.Content.Find.Execute("212×3=636", $true, $false, $false, $false, $false, $true, 1, $false, "418×5=2090", 2) | Out-Null
$d.Content.Find.Execute("662×7=4634", $true, $false, $false, $false, $false, $true, 1, $false, "539×5=2695", 2) | Out-Null
$d.Content.Find.Execute("928×5=4640", $true, $false, $false, $false, $false, $true, 1, $false, "476×9=4284", 2) | Out-Null
$d.Content.Find.Execute("190×6=1140", $true, $false, $false, $false, $false, $true, 1, $false, "906×5=4530", 2) | Out-Null
$d.Content.Find.Execute("558×9=5022", $true, $false, $false, $false, $false, $true, 1, $false, "510×3=1530", 2) | Out-Null
$d.Content.Find.Execute("863×4=3452", $true, $false, $false, $false, $false, $true, 1, $false, "314×2=628", 2) | Out-Null
$d.Content.Find.Execute("510×2=1020", $true, $false, $false, $false, $false, $true, 1, $false, "773×7=5411", 2) | Out-Null
$d.Content.Find.Execute("831×4=3324", $true, $false, $false, $false, $false, $true, 1, $false, "405×6=2430", 2) | Out-Null
